{"js": "// Resume update: rewrite the \"Projects include\" bullet block under the\n// lululemon job into five new bullets (one summary line + four details),\n// promote the four detail bullets from sub-level to top-level, and move the\n// \"_GoBack\" bookmark to the very start of the document (matches a resave by\n// Word after the content edit).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Exact paragraph-text replacements (old -> new). Order follows the\n// document; each is matched against the current paragraph text so the\n// script is resilient to the exact paragraph index.\nconst replacements = [\n  {\n    match: \"Projects include:\",\n    text: \"Developed Omnichannel middleware API services (Java) and web dashboard business tools (Angular JS)\",\n    promote: false,\n  },\n  {\n    match: \"Omnichannel API services and dashboard business tools\",\n    text: \"Created an iOS application feature to allow streamlining in store processes\",\n    promote: true,\n  },\n  {\n    match: \"iOS application feature to allow streamlining in store processes\",\n    text: \"Built a second iOS application to integrate with third party applications for mobile printing to help store employees process orders (Objective-C)\",\n    promote: true,\n  },\n  {\n    match: \"Website customer facing redesign of certain features\",\n    text: \"Redesign customer facing website features\",\n    promote: true,\n  },\n  {\n    match: \"Developed a API test framework to automate QA test cases and integrate in build processes\",\n    text: \"Developed an API test framework to automate QA test cases to integrate with build processes\",\n    promote: true,\n  },\n];\n\nconst items = paragraphs.items;\nfor (const rep of replacements) {\n  const para = items.find((p) => p.text === rep.match);\n  if (!para) {\n    continue;\n  }\n  para.insertText(rep.text, Word.InsertLocation.replace);\n  if (rep.promote) {\n    // These bullets move from the nested level (ilvl=1, indent 1260\n    // twips = 63pt) up to the top level (ilvl=0, indent 360 twips = 18pt).\n    para.listItemOrNullObject.level = 0;\n    para.leftIndent = 18;\n  }\n}\nawait context.sync();\n\n// Relocate the \"_GoBack\" bookmark (an artifact Word writes at the most\n// recent edit location) from the end of the document to the very start.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\nbody.getRange(Word.RangeLocation.start).insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Resume update: rewrite the \"Projects include\" bullet block under the\n# lululemon job into five new bullets (one summary line + four details),\n# promote the four detail bullets from sub-level to top-level, and move the\n# \"_GoBack\" bookmark to the very start of the document (matches a resave by\n# Word after the content edit).\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($doc, $text) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs($i)\n        $t = $p.Range.Text\n        if ($t.Length -gt 0) {\n            $t = $t.Substring(0, $t.Length - 1)\n        }\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\n$replacements = @(\n    @{ Match = \"Projects include:\"; Text = \"Developed Omnichannel middleware API services (Java) and web dashboard business tools (Angular JS)\"; Promote = $false },\n    @{ Match = \"Omnichannel API services and dashboard business tools\"; Text = \"Created an iOS application feature to allow streamlining in store processes\"; Promote = $true },\n    @{ Match = \"iOS application feature to allow streamlining in store processes\"; Text = \"Built a second iOS application to integrate with third party applications for mobile printing to help store employees process orders (Objective-C)\"; Promote = $true },\n    @{ Match = \"Website customer facing redesign of certain features\"; Text = \"Redesign customer facing website features\"; Promote = $true },\n    @{ Match = \"Developed a API test framework to automate QA test cases and integrate in build processes\"; Text = \"Developed an API test framework to automate QA test cases to integrate with build processes\"; Promote = $true }\n)\n\nforeach ($rep in $replacements) {\n    $p = Find-ParagraphByText $d $rep.Match\n    if ($p -ne $null) {\n        $p.Range.Text = $rep.Text\n        if ($rep.Promote) {\n            # These bullets move from the nested level (ilvl=1, indent 1260\n            # twips = 63pt) up to the top level (ilvl=0, indent 360 twips = 18pt).\n            $p.Range.ListFormat.ListLevelNumber = 1\n            $p.LeftIndent = 18\n        }\n    }\n}\n\n# Relocate the \"_GoBack\" bookmark (an artifact Word writes at the most\n# recent edit location) from the end of the document to the very start.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n# Word's COM bookmark-add has an edge case at absolute position 0, so we\n# insert a temporary marker character, bookmark right after it, then\n# remove the marker -- leaving a zero-length bookmark at the true start.\n$startRange = $d.Range(0, 0)\n$startRange.InsertBefore(\"X\")\n$afterMarker = $d.Range(1, 1)\n$d.Bookmarks.Add(\"_GoBack\", $afterMarker)\n$d.Range(0, 1).Delete()\n"}
